# Generate Report for Handoff
# Adds two new tracked files (5a12b1f2-...md and 793a859c-...md) to each of
# the three report sheets (Overview, zh-cn, de-de), extending their tables
# and hyperlinks accordingly.

$wb = $excel.ActiveWorkbook

$dateStamp = "2016-09-04 22:45:07"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 4 - 5a12b1f2-...md
$wsOverview.Cells.Item(4,1).Value2 = "5a12b1f2-57b0-4e1d-97a2-125898b65dbe.md"
$wsOverview.Cells.Item(4,2).Value2 = "e2e\5a12b1f2-57b0-4e1d-97a2-125898b65dbe.md"
$wsOverview.Cells.Item(4,3).Value2 = ".md"
$wsOverview.Cells.Item(4,5).Value2 = "Ready for handoff"
$wsOverview.Cells.Item(4,6).Value2 = "Ready for handoff"
$wsOverview.Cells.Item(4,7).Value2 = $dateStamp
$wsOverview.Cells.Item(4,7).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8bc95a6c29194b430143ab3e463dcc7c6ab23edb/e2e/5a12b1f2-57b0-4e1d-97a2-125898b65dbe.md", "", "", "e2e\5a12b1f2-57b0-4e1d-97a2-125898b65dbe.md")

# Row 5 - 793a859c-...md
$wsOverview.Cells.Item(5,1).Value2 = "793a859c-b0dd-4717-8005-929d2c061933.md"
$wsOverview.Cells.Item(5,2).Value2 = "e2e\793a859c-b0dd-4717-8005-929d2c061933.md"
$wsOverview.Cells.Item(5,3).Value2 = ".md"
$wsOverview.Cells.Item(5,5).Value2 = "Ready for handoff"
$wsOverview.Cells.Item(5,6).Value2 = "Ready for handoff"
$wsOverview.Cells.Item(5,7).Value2 = $dateStamp
$wsOverview.Cells.Item(5,7).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/75985593933061f69a9256f73eb1e92c7f89e07b/e2e/793a859c-b0dd-4717-8005-929d2c061933.md", "", "", "e2e\793a859c-b0dd-4717-8005-929d2c061933.md")

$wsOverview.Range("B4").Style = "HyperLink"
$wsOverview.Range("B5").Style = "HyperLink"

# Grow the "Overview" table (table3.xml) to include the two new rows.
$wsOverview.ListObjects.Item(1).Resize($wsOverview.Range("A1:G5"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Row 4 - 5a12b1f2-...md
$wsZhCn.Cells.Item(4,1).Value2 = "5a12b1f2-57b0-4e1d-97a2-125898b65dbe.md"
$wsZhCn.Cells.Item(4,2).Value2 = ".md"
$wsZhCn.Cells.Item(4,3).Value2 = "Ready for handoff"
$wsZhCn.Cells.Item(4,4).Value2 = "e2e"
$wsZhCn.Cells.Item(4,5).Value2 = "ht"
$wsZhCn.Cells.Item(4,6).Value2 = "'False"
$wsZhCn.Cells.Item(4,7).Value2 = "5a12b1f2-57b0-4e1d-97a2-125898b65dbe.d965337e8693407b4a8b61dc30aa9c43e39cc158.zh-cn.xlf"
$wsZhCn.Cells.Item(4,8).Value2 = "2016-09-04 22:44:58"
$wsZhCn.Cells.Item(4,8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Cells.Item(4,11).Value2 = "0001-01-01 00:00:00"
$wsZhCn.Cells.Item(4,11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Cells.Item(4,13).Value2 = "'True"
$wsZhCn.Cells.Item(4,15).Value2 = "'False"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/1b289a8bb4dd9930eb765d3673389d1736a59ff5/e2e/5a12b1f2-57b0-4e1d-97a2-125898b65dbe.md", "", "", "5a12b1f2-57b0-4e1d-97a2-125898b65dbe.md")

# Row 5 - 793a859c-...md
$wsZhCn.Cells.Item(5,1).Value2 = "793a859c-b0dd-4717-8005-929d2c061933.md"
$wsZhCn.Cells.Item(5,2).Value2 = ".md"
$wsZhCn.Cells.Item(5,3).Value2 = "Ready for handoff"
$wsZhCn.Cells.Item(5,4).Value2 = "e2e"
$wsZhCn.Cells.Item(5,5).Value2 = "ht"
$wsZhCn.Cells.Item(5,6).Value2 = "'False"
$wsZhCn.Cells.Item(5,7).Value2 = "793a859c-b0dd-4717-8005-929d2c061933.02e66577b5d6907ceace537444ef3fea521468b3.zh-cn.xlf"
$wsZhCn.Cells.Item(5,8).Value2 = "2016-09-04 22:44:58"
$wsZhCn.Cells.Item(5,8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Cells.Item(5,11).Value2 = "0001-01-01 00:00:00"
$wsZhCn.Cells.Item(5,11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Cells.Item(5,13).Value2 = "'True"
$wsZhCn.Cells.Item(5,15).Value2 = "'False"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/1b289a8bb4dd9930eb765d3673389d1736a59ff5/e2e/793a859c-b0dd-4717-8005-929d2c061933.md", "", "", "793a859c-b0dd-4717-8005-929d2c061933.md")

# Grow the "zh-cn" table (table1.xml) to include the two new rows.
$wsZhCn.ListObjects.Item(1).Resize($wsZhCn.Range("A1:P5"))

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 4 - 5a12b1f2-...md
$wsDeDe.Cells.Item(4,1).Value2 = "5a12b1f2-57b0-4e1d-97a2-125898b65dbe.md"
$wsDeDe.Cells.Item(4,2).Value2 = ".md"
$wsDeDe.Cells.Item(4,3).Value2 = "Ready for handoff"
$wsDeDe.Cells.Item(4,4).Value2 = "e2e"
$wsDeDe.Cells.Item(4,5).Value2 = "ht"
$wsDeDe.Cells.Item(4,6).Value2 = "'False"
$wsDeDe.Cells.Item(4,7).Value2 = "5a12b1f2-57b0-4e1d-97a2-125898b65dbe.d965337e8693407b4a8b61dc30aa9c43e39cc158.de-de.xlf"
$wsDeDe.Cells.Item(4,8).Value2 = $dateStamp
$wsDeDe.Cells.Item(4,8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Cells.Item(4,11).Value2 = "0001-01-01 00:00:00"
$wsDeDe.Cells.Item(4,11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Cells.Item(4,13).Value2 = "'True"
$wsDeDe.Cells.Item(4,15).Value2 = "'False"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/4abac84ce6b65d220b7f640998db3a8f2a79ec19/e2e/5a12b1f2-57b0-4e1d-97a2-125898b65dbe.md", "", "", "5a12b1f2-57b0-4e1d-97a2-125898b65dbe.md")

# Row 5 - 793a859c-...md
$wsDeDe.Cells.Item(5,1).Value2 = "793a859c-b0dd-4717-8005-929d2c061933.md"
$wsDeDe.Cells.Item(5,2).Value2 = ".md"
$wsDeDe.Cells.Item(5,3).Value2 = "Ready for handoff"
$wsDeDe.Cells.Item(5,4).Value2 = "e2e"
$wsDeDe.Cells.Item(5,5).Value2 = "ht"
$wsDeDe.Cells.Item(5,6).Value2 = "'False"
$wsDeDe.Cells.Item(5,7).Value2 = "793a859c-b0dd-4717-8005-929d2c061933.02e66577b5d6907ceace537444ef3fea521468b3.de-de.xlf"
$wsDeDe.Cells.Item(5,8).Value2 = $dateStamp
$wsDeDe.Cells.Item(5,8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Cells.Item(5,11).Value2 = "0001-01-01 00:00:00"
$wsDeDe.Cells.Item(5,11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Cells.Item(5,13).Value2 = "'True"
$wsDeDe.Cells.Item(5,15).Value2 = "'False"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/4abac84ce6b65d220b7f640998db3a8f2a79ec19/e2e/793a859c-b0dd-4717-8005-929d2c061933.md", "", "", "793a859c-b0dd-4717-8005-929d2c061933.md")

# Grow the "de-de" table (table2.xml) to include the two new rows.
$wsDeDe.ListObjects.Item(1).Resize($wsDeDe.Range("A1:P5"))
